$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F7").Value = '135_product_information'
$ws.Range("F8").Value = '135_product_information'
$ws.Range("F9").Value = '135_product_information'
$ws.Range("F24").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F25").Value = 'ppe'
$ws.Range("F26").Value = 'ppe'
$ws.Range("F27").Value = 'ppe'
$ws.Range("F30").Value = 'off target movement || application instructions || env warning - species || env warning - water'
$ws.Range("F34").Value = 'application instructions'
$ws.Range("F35").Value = '134_non-agriculture_use_requirements'
$ws.Range("F36").Value = '134_non-agriculture_use_requirements'
$ws.Range("F38").Value = '93_referral_statement || chemigation'
$ws.Range("F45").Value = 'off target movement'
$ws.Range("F46").Value = '172_sensitive_areas || off target movement'
$ws.Range("F47").Value = 'application instructions'
$ws.Range("F48").Value = 'mixing'
$ws.Range("F49").Value = 'mixing'
$ws.Range("F50").Value = 'mixing'
$ws.Range("F52").Value = 'use restrictions'
$ws.Range("F68").Value = 'application instructions'
$ws.Range("F72").Value = 'application instructions'
$ws.Range("F73").Value = 'application instructions'
$ws.Range("F74").Value = 'use restrictions'
$ws.Range("F90").Value = 'application instructions'
$ws.Range("F91").Value = 'application instructions'
$ws.Range("F92").Value = 'application instructions'
$ws.Range("F109").Value = 'application instructions'
$ws.Range("F128").Value = 'application instructions'
$ws.Range("F129").Value = 'use restrictions'
$ws.Range("F143").Value = 'application instructions'
$ws.Range("F192").Value = 'application instructions'
$ws.Range("F212").Value = 'application instructions'
$ws.Range("F213").Value = 'application instructions'
$ws.Range("F214").Value = 'application instructions'
$ws.Range("F216").Value = 'use restrictions'
$ws.Range("F240").Value = 'application instructions'
$ws.Range("F280").Value = 'application instructions'
$ws.Range("F298").Value = 'application instructions'
$ws.Range("F313").Value = 'application instructions'
$ws.Range("F329").Value = 'application instructions'
$ws.Range("F349").Value = 'application instructions'
$ws.Range("F350").Value = 'application instructions'
$ws.Range("F353").Value = 'application instructions'
$ws.Range("F355").Value = 'mixing || application instructions'
$ws.Range("F358").Value = 'use restrictions'
$ws.Range("F379").Value = 'application instructions'
$ws.Range("F399").Value = 'application instructions'
$ws.Range("F400").Value = 'use restrictions'
$ws.Range("F413").Value = 'application instructions'
$ws.Range("F448").Value = 'application instructions'
$ws.Range("F449").Value = 'application instructions'
$ws.Range("F450").Value = 'application instructions'
$ws.Range("F451").Value = 'application instructions'
$ws.Range("F455").Value = 'application instructions'
$ws.Range("F467").Value = 'use restrictions'
$ws.Range("F469").Value = 'use restrictions'
$ws.Range("F470").Value = 'application instructions'
$ws.Range("F471").Value = 'application instructions'
$ws.Range("F511").Value = 'use restrictions'
$ws.Range("F512").Value = 'application instructions'
$ws.Range("F518").Value = '172_sensitive_areas || off target movement'
$ws.Range("F532").Value = '154_pesticide_storage'
